# Apply the PO data update described in the commit:
#  1. Rename header labels on the existing sheets
#  2. Add a new "PO Forecast" sheet with forecast data

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy the header/date cell formatting from the "Weekly Quantity" sheet so the
# new sheet reuses the same styles (bold/centered header, date number format)
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row
$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

# Forecast data rows: ds, PO_Forecast, yhat_lower, yhat_upper
$data = @(
    @(45515.99999999999, 86, 17.23451431900819, 155.0078979528322),
    @(45578.99999999999, 70, 7.877701468294379, 138.6486072160009),
    @(45627.99999999999, 57, -11.24188739851971, 122.9054990346986),
    @(45634.99999999999, 55, -9.059238446172566, 120.80453882575),
    @(45641.99999999999, 53, -9.446251332856182, 117.7421981264949),
    @(45648.99999999999, 52, -11.49175332626255, 113.4399891786674),
    @(45655.99999999999, 50, -19.09381489943021, 113.7289228965717),
    @(45662.99999999999, 48, -14.11894413183844, 122.0481631063747),
    @(45669.99999999999, 46, -18.84648091560007, 114.7148244636254),
    @(45676.99999999999, 44, -28.52159797220614, 107.0298892986481),
    @(45683.99999999999, 43, -19.28658335867257, 110.0053673704274),
    @(45690.99999999999, 41, -23.68010210631523, 107.6093226651386)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

Write-Host "PO Forecast sheet added and headers updated"
